$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheetId 1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 3910
$ws1.Range("G5").Value = 59
$ws1.Range("F7").Value = 2576
$ws1.Range("F9").Value = 3184
$ws1.Range("F11").Value = 2333
$ws1.Range("F15").Value = 467
$ws1.Range("F20").Value = 312
$ws1.Range("F21").Value = 434
$ws1.Range("F22").Value = 673
$ws1.Range("F23").Value = 1418
$ws1.Range("F27").Value = 140
$ws1.Range("F28").Value = 158
$ws1.Range("F29").Value = 37
$ws1.Range("F31").Value = 69
$ws1.Range("F32").Value = 4379
$ws1.Range("F33").Value = 4218
$ws1.Range("F34").Value = 86
$ws1.Range("F35").Value = 138
$ws1.Range("F38").Value = 1151
$ws1.Range("F43").Value = 181
$ws1.Range("F45").Value = 112
$ws1.Range("F46").Value = 43
$ws1.Range("F48").Value = 65

# Sheet "演出" (sheetId 2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F9").Value = 17
$ws2.Range("F16").Value = 214

# Sheet "本地生活" (sheetId 3)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 156
$ws3.Range("F4").Value = 2330

# Sheet "全部类型" (sheetId 4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 156
$ws4.Range("F10").Value = 3910
$ws4.Range("G10").Value = 59
$ws4.Range("F12").Value = 2576
$ws4.Range("F14").Value = 3184
$ws4.Range("F17").Value = 2333
$ws4.Range("F24").Value = 312
$ws4.Range("F25").Value = 434
$ws4.Range("F26").Value = 673
$ws4.Range("F27").Value = 1418
$ws4.Range("F30").Value = 158
$ws4.Range("F32").Value = 69
$ws4.Range("F33").Value = 17
$ws4.Range("F34").Value = 4379
$ws4.Range("F35").Value = 4218
$ws4.Range("F36").Value = 86
$ws4.Range("F38").Value = 1151
$ws4.Range("F46").Value = 181
$ws4.Range("F47").Value = 112
$ws4.Range("F48").Value = 43
$ws4.Range("F49").Value = 214
